$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text formatting instead of
# being auto-converted to numbers by Excel when values look numeric.
$ws.Range("D2:D51").NumberFormat = "@"

$updates = @(
    @{ Addr = 'D2'; Val = '29.723.97' }
    @{ Addr = 'E2'; Val = '  +0.55%  ' }
    @{ Addr = 'D3'; Val = '1.851.48' }
    @{ Addr = 'E3'; Val = '  +0.44%  ' }
    @{ Addr = 'D4'; Val = '0.9997' }
    @{ Addr = 'E4'; Val = '  +0.10%  ' }
    @{ Addr = 'D5'; Val = '243.67' }
    @{ Addr = 'E5'; Val = '  -0.33%  ' }
    @{ Addr = 'D6'; Val = '0.6568' }
    @{ Addr = 'E6'; Val = '  +4.04%  ' }
    @{ Addr = 'D7'; Val = '1.000' }
    @{ Addr = 'E7'; Val = '  +0.08%  ' }
    @{ Addr = 'B8'; Val = 'Dogecoin' }
    @{ Addr = 'C8'; Val = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge' }
    @{ Addr = 'D8'; Val = '0.07498' }
    @{ Addr = 'E8'; Val = '  +0.71%  ' }
    @{ Addr = 'B9'; Val = 'Cardano' }
    @{ Addr = 'C9'; Val = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada' }
    @{ Addr = 'D9'; Val = '0.2981' }
    @{ Addr = 'E9'; Val = '  +1.29%  ' }
    @{ Addr = 'B10'; Val = 'Solana' }
    @{ Addr = 'C10'; Val = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol' }
    @{ Addr = 'D10'; Val = '24.54' }
    @{ Addr = 'E10'; Val = '  +2.86%  ' }
    @{ Addr = 'B11'; Val = 'TRON' }
    @{ Addr = 'C11'; Val = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx' }
    @{ Addr = 'D11'; Val = '0.07648' }
    @{ Addr = 'E11'; Val = '  -0.36%  ' }
    @{ Addr = 'B12'; Val = 'WrappedEther' }
    @{ Addr = 'C12'; Val = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' }
    @{ Addr = 'D12'; Val = '1.852.94' }
    @{ Addr = 'E12'; Val = '  +0.64%  ' }
    @{ Addr = 'B13'; Val = 'Polkadot' }
    @{ Addr = 'C13'; Val = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot' }
    @{ Addr = 'D13'; Val = '5.051' }
    @{ Addr = 'E13'; Val = '  +0.53%  ' }
    @{ Addr = 'B14'; Val = 'Polygon' }
    @{ Addr = 'C14'; Val = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic' }
    @{ Addr = 'D14'; Val = '0.6871' }
    @{ Addr = 'E14'; Val = '  +0.86%  ' }
    @{ Addr = 'B15'; Val = 'Litecoin' }
    @{ Addr = 'C15'; Val = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc' }
    @{ Addr = 'D15'; Val = '83.56' }
    @{ Addr = 'E15'; Val = '  -0.72%  ' }
    @{ Addr = 'B16'; Val = 'ShibaInu' }
    @{ Addr = 'C16'; Val = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib' }
    @{ Addr = 'D16'; Val = '0.000009602' }
    @{ Addr = 'E16'; Val = '  +3.02%  ' }
    @{ Addr = 'B17'; Val = 'Uniswap' }
    @{ Addr = 'C17'; Val = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni' }
    @{ Addr = 'D17'; Val = '6.132' }
    @{ Addr = 'E17'; Val = '  +2.63%  ' }
    @{ Addr = 'B18'; Val = 'WrappedBTC' }
    @{ Addr = 'C18'; Val = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc' }
    @{ Addr = 'D18'; Val = '29.755.16' }
    @{ Addr = 'E18'; Val = '  +0.77%  ' }
    @{ Addr = 'B19'; Val = 'WrappedliquidstakedEther2.0' }
    @{ Addr = 'C19'; Val = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth' }
    @{ Addr = 'D19'; Val = '2.108.20' }
    @{ Addr = 'E19'; Val = '  +0.79%  ' }
    @{ Addr = 'B20'; Val = 'BitcoinCash' }
    @{ Addr = 'C20'; Val = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch' }
    @{ Addr = 'D20'; Val = '237.41' }
    @{ Addr = 'E20'; Val = '  -0.21%  ' }
    @{ Addr = 'B21'; Val = 'Avalanche' }
    @{ Addr = 'C21'; Val = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax' }
    @{ Addr = 'D21'; Val = '12.62' }
    @{ Addr = 'E21'; Val = '  +0.30%  ' }
    @{ Addr = 'B22'; Val = 'Dai' }
    @{ Addr = 'C22'; Val = 'https://coinranking.com/coin/MoTuySvg7+dai-dai' }
    @{ Addr = 'D22'; Val = '1.000' }
    @{ Addr = 'E22'; Val = '  +0.06%  ' }
    @{ Addr = 'B23'; Val = 'Chainlink' }
    @{ Addr = 'C23'; Val = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link' }
    @{ Addr = 'D23'; Val = '7.726' }
    @{ Addr = 'E23'; Val = '  +4.93%  ' }
    @{ Addr = 'B24'; Val = 'BinanceUSD' }
    @{ Addr = 'C24'; Val = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd' }
    @{ Addr = 'D24'; Val = '1.001' }
    @{ Addr = 'E24'; Val = '  +0.18%  ' }
    @{ Addr = 'B25'; Val = 'Monero' }
    @{ Addr = 'C25'; Val = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' }
    @{ Addr = 'D25'; Val = '158.29' }
    @{ Addr = 'E25'; Val = '  -0.28%  ' }
    @{ Addr = 'B26'; Val = 'Stellar' }
    @{ Addr = 'C26'; Val = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm' }
    @{ Addr = 'D26'; Val = '0.1427' }
    @{ Addr = 'E26'; Val = '  +0.48%  ' }
    @{ Addr = 'B27'; Val = 'Cosmos' }
    @{ Addr = 'C27'; Val = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' }
    @{ Addr = 'D27'; Val = '8.545' }
    @{ Addr = 'E27'; Val = '  +0.27%  ' }
    @{ Addr = 'B28'; Val = 'EthereumClassic' }
    @{ Addr = 'C28'; Val = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc' }
    @{ Addr = 'D28'; Val = '17.86' }
    @{ Addr = 'E28'; Val = '  +0.25%  ' }
    @{ Addr = 'B29'; Val = 'PancakeSwap' }
    @{ Addr = 'C29'; Val = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake' }
    @{ Addr = 'D29'; Val = '1.494' }
    @{ Addr = 'E29'; Val = '  -0.19%  ' }
    @{ Addr = 'D30'; Val = '0.06046' }
    @{ Addr = 'E30'; Val = '  +0.47%  ' }
    @{ Addr = 'B31'; Val = 'Toncoin' }
    @{ Addr = 'C31'; Val = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton' }
    @{ Addr = 'D31'; Val = '1.273' }
    @{ Addr = 'E31'; Val = '  +1.57%  ' }
    @{ Addr = 'B32'; Val = 'Filecoin' }
    @{ Addr = 'C32'; Val = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' }
    @{ Addr = 'D32'; Val = '4.143' }
    @{ Addr = 'E32'; Val = '  +0.61%  ' }
    @{ Addr = 'B33'; Val = 'InternetComputer(DFINITY)' }
    @{ Addr = 'C33'; Val = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' }
    @{ Addr = 'D33'; Val = '4.070' }
    @{ Addr = 'E33'; Val = '  -1.95%  ' }
    @{ Addr = 'B34'; Val = 'LidoDAOToken' }
    @{ Addr = 'C34'; Val = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo' }
    @{ Addr = 'D34'; Val = '1.871' }
    @{ Addr = 'E34'; Val = '  +0.19%  ' }
    @{ Addr = 'B35'; Val = 'ARBITRUM' }
    @{ Addr = 'C35'; Val = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb' }
    @{ Addr = 'D35'; Val = '1.182' }
    @{ Addr = 'E35'; Val = '  +2.94%  ' }
    @{ Addr = 'B36'; Val = 'ImmutableX' }
    @{ Addr = 'C36'; Val = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' }
    @{ Addr = 'D36'; Val = '0.7269' }
    @{ Addr = 'E36'; Val = '  +0.17%  ' }
    @{ Addr = 'B37'; Val = 'HuobiToken' }
    @{ Addr = 'C37'; Val = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht' }
    @{ Addr = 'D37'; Val = '2.600' }
    @{ Addr = 'E37'; Val = '  -0.35%  ' }
    @{ Addr = 'B38'; Val = 'MXToken' }
    @{ Addr = 'C38'; Val = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' }
    @{ Addr = 'D38'; Val = '2.803' }
    @{ Addr = 'E38'; Val = '  -2.73%  ' }
    @{ Addr = 'B39'; Val = 'VeChain' }
    @{ Addr = 'C39'; Val = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' }
    @{ Addr = 'D39'; Val = '0.01791' }
    @{ Addr = 'E39'; Val = '  +1.23%  ' }
    @{ Addr = 'B40'; Val = 'Maker' }
    @{ Addr = 'C40'; Val = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr' }
    @{ Addr = 'D40'; Val = '1.203.44' }
    @{ Addr = 'E40'; Val = '  -1.53%  ' }
    @{ Addr = 'B41'; Val = 'FraxShare' }
    @{ Addr = 'C41'; Val = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' }
    @{ Addr = 'D41'; Val = '6.289' }
    @{ Addr = 'E41'; Val = '  -0.05%  ' }
    @{ Addr = 'B42'; Val = 'TrustWalletToken' }
    @{ Addr = 'C42'; Val = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' }
    @{ Addr = 'D42'; Val = '0.9134' }
    @{ Addr = 'E42'; Val = '  -1.23%  ' }
    @{ Addr = 'B43'; Val = 'PaxDollar' }
    @{ Addr = 'C43'; Val = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp' }
    @{ Addr = 'D43'; Val = '0.9999' }
    @{ Addr = 'E43'; Val = '  -0.11%  ' }
    @{ Addr = 'B44'; Val = 'RocketPoolETH' }
    @{ Addr = 'C44'; Val = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth' }
    @{ Addr = 'D44'; Val = '2.015.73' }
    @{ Addr = 'E44'; Val = '  +0.34%  ' }
    @{ Addr = 'B45'; Val = 'Quant' }
    @{ Addr = 'C45'; Val = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt' }
    @{ Addr = 'D45'; Val = '101.09' }
    @{ Addr = 'E45'; Val = '  -0.77%  ' }
    @{ Addr = 'B46'; Val = 'Aave' }
    @{ Addr = 'C46'; Val = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' }
    @{ Addr = 'D46'; Val = '66.63' }
    @{ Addr = 'E46'; Val = '  +0.99%  ' }
    @{ Addr = 'B47'; Val = 'Aptos' }
    @{ Addr = 'C47'; Val = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' }
    @{ Addr = 'D47'; Val = '7.332' }
    @{ Addr = 'E47'; Val = '  +9.84%  ' }
    @{ Addr = 'D48'; Val = '0.00000000118' }
    @{ Addr = 'E48'; Val = '  -1.83%  ' }
    @{ Addr = 'B49'; Val = 'TheSandbox' }
    @{ Addr = 'C49'; Val = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand' }
    @{ Addr = 'D49'; Val = '0.4057' }
    @{ Addr = 'E49'; Val = '  -0.35%  ' }
    @{ Addr = 'B50'; Val = 'EnergySwap' }
    @{ Addr = 'C50'; Val = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' }
    @{ Addr = 'D50'; Val = '9.126' }
    @{ Addr = 'E50'; Val = '  -1.51%  ' }
    @{ Addr = 'B51'; Val = 'RenderToken' }
    @{ Addr = 'C51'; Val = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' }
    @{ Addr = 'D51'; Val = '1.671' }
    @{ Addr = 'E51'; Val = '  +3.57%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Addr).Value = $u.Val
}
